$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restore the (previously minimized) window to a normal state ---
$wb.Windows.Item(1).WindowState = -4143

# --- Row 14: re-purpose the "Organstatus evaluering kirurgi" entry.
#     B14 now holds "Undersøkelse", D14 now holds the text that used to be
#     in D14's neighbour cell ("Skal ikke vises fra start, kun mulighet ...") ---
$ws.Cells.Item(14, 2).Value = "Undersøkelse"
$ws.Cells.Item(14, 4).Value = "Skal ikke vises fra start, kun mulighet å legge til (som diagnose)"

# --- New row 15: Auskultasjon thorax ---
$ws.Cells.Item(15, 2).Value = "Auskultasjon thorax"
$ws.Cells.Item(15, 3).Value = "Alle elementer som er brukt"
$ws.Cells.Item(15, 4).Value = "Skal ikke vises fra start, mulighet til å legge til og velge f.eks. legg til asukultasjon hjerte, legg til auskultasjon lunge etc. Velge fra en nedtrekksmeny?"
$ws.Cells.Item(15, 5).Value = "BNA/MGR"

# Copy formatting (wrap text, style) from the row above, and give the new
# "Kommentar" cell (E15) the same yellow highlight style used on E14/E9.
$ws.Cells.Item(14, 2).Copy()
$ws.Cells.Item(15, 2).PasteSpecial(-4122)
$ws.Cells.Item(14, 3).Copy()
$ws.Cells.Item(15, 3).PasteSpecial(-4122)
$ws.Cells.Item(14, 4).Copy()
$ws.Cells.Item(15, 4).PasteSpecial(-4122)
$ws.Cells.Item(14, 5).Copy()
$ws.Cells.Item(15, 5).PasteSpecial(-4122)

# Restore the values (PasteSpecial of formats only shouldn't overwrite them,
# but set again defensively in case paste touched content)
$ws.Cells.Item(15, 2).Value = "Auskultasjon thorax"
$ws.Cells.Item(15, 3).Value = "Alle elementer som er brukt"
$ws.Cells.Item(15, 4).Value = "Skal ikke vises fra start, mulighet til å legge til og velge f.eks. legg til asukultasjon hjerte, legg til auskultasjon lunge etc. Velge fra en nedtrekksmeny?"
$ws.Cells.Item(15, 5).Value = "BNA/MGR"

# Row 15 wraps onto several lines, same as row 4's long comment.
$ws.Rows.Item(15).RowHeight = 43.2

# --- Update the active selection to match the new end of the sheet ---
$null = $ws.Range("C16").Select()
